# Fix ecdc-csv data for 2020-03-29.
# Appends a new row (row 91) for the date 2020-03-29 to both the
# "Confirmed" and "Deaths" worksheets, with per-country case/death counts.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Column A (date, stored as shared string "2020-03-29") ---
# Use a string-literal formula then paste-as-values so the final cell is a
# plain shared string (t="s") rather than triggering Excel's automatic
# text-to-date conversion (which would turn it into a date serial number)
# and without leaving stray style/number-format entries behind.
$ws1.Cells.Item(91, 1).Formula = "=""2020-03-29"""
$ws1.Cells.Item(91, 1).Copy()
$ws1.Cells.Item(91, 1).PasteSpecial(-4163)

$ws2.Cells.Item(91, 1).Formula = "=""2020-03-29"""
$ws2.Cells.Item(91, 1).Copy()
$ws2.Cells.Item(91, 1).PasteSpecial(-4163)

$excel.CutCopyMode = 0

# --- Confirmed sheet (sheet1), row 91 ---
$ws1.Cells.Item(91, 2).Value = 15
$ws1.Cells.Item(91, 3).Value = 11
$ws1.Cells.Item(91, 4).Value = 104
$ws1.Cells.Item(91, 5).Value = 41
$ws1.Cells.Item(91, 6).Value = 0
$ws1.Cells.Item(91, 7).Value = 0
$ws1.Cells.Item(91, 8).Value = 0
$ws1.Cells.Item(91, 9).Value = 55
$ws1.Cells.Item(91, 10).Value = 52
$ws1.Cells.Item(91, 11).Value = 0
$ws1.Cells.Item(91, 12).Value = 431
$ws1.Cells.Item(91, 13).Value = 594
$ws1.Cells.Item(91, 14).Value = 17
$ws1.Cells.Item(91, 15).Value = 2
$ws1.Cells.Item(91, 16).Value = 7
$ws1.Cells.Item(91, 17).Value = 0
$ws1.Cells.Item(91, 18).Value = 2
$ws1.Cells.Item(91, 19).Value = 0
$ws1.Cells.Item(91, 20).Value = 1850
$ws1.Cells.Item(91, 21).Value = 0
$ws1.Cells.Item(91, 22).Value = 0
$ws1.Cells.Item(91, 23).Value = 5
$ws1.Cells.Item(91, 24).Value = 1
$ws1.Cells.Item(91, 25).Value = 7
$ws1.Cells.Item(91, 26).Value = 28
$ws1.Cells.Item(91, 27).Value = 487
$ws1.Cells.Item(91, 28).Value = 0
$ws1.Cells.Item(91, 29).Value = 5
$ws1.Cells.Item(91, 30).Value = 38
$ws1.Cells.Item(91, 31).Value = 34
$ws1.Cells.Item(91, 32).Value = 1
$ws1.Cells.Item(91, 33).Value = 11
$ws1.Cells.Item(91, 34).Value = 697
$ws1.Cells.Item(91, 35).Value = 1
$ws1.Cells.Item(91, 37).Value = 0
$ws1.Cells.Item(91, 38).Value = 1
$ws1.Cells.Item(91, 39).Value = 0
$ws1.Cells.Item(91, 40).Value = 299
$ws1.Cells.Item(91, 41).Value = 129
$ws1.Cells.Item(91, 42).Value = 69
$ws1.Cells.Item(91, 43).Value = 15
$ws1.Cells.Item(91, 44).Value = 32
$ws1.Cells.Item(91, 45).Value = 39
$ws1.Cells.Item(91, 46).Value = 71
$ws1.Cells.Item(91, 47).Value = 39
$ws1.Cells.Item(91, 48).Value = 0
$ws1.Cells.Item(91, 49).Value = 17
$ws1.Cells.Item(91, 50).Value = 384
$ws1.Cells.Item(91, 51).Value = 0
$ws1.Cells.Item(91, 52).Value = 155
$ws1.Cells.Item(91, 53).Value = 2
$ws1.Cells.Item(91, 54).Value = 0
$ws1.Cells.Item(91, 55).Value = 138
$ws1.Cells.Item(91, 56).Value = 208
$ws1.Cells.Item(91, 57).Value = 41
$ws1.Cells.Item(91, 58).Value = 6
$ws1.Cells.Item(91, 59).Value = 0
$ws1.Cells.Item(91, 60).Value = 0
$ws1.Cells.Item(91, 61).Value = 65
$ws1.Cells.Item(91, 62).Value = 0
$ws1.Cells.Item(91, 63).Value = 0
$ws1.Cells.Item(91, 64).Value = 11
$ws1.Cells.Item(91, 65).Value = 0
$ws1.Cells.Item(91, 66).Value = 193
$ws1.Cells.Item(91, 67).Value = 4611
$ws1.Cells.Item(91, 68).Value = 4
$ws1.Cells.Item(91, 69).Value = 0
$ws1.Cells.Item(91, 70).Value = 0
$ws1.Cells.Item(91, 71).Value = 4
$ws1.Cells.Item(91, 72).Value = 3965
$ws1.Cells.Item(91, 73).Value = 4
$ws1.Cells.Item(91, 74).Value = 1
$ws1.Cells.Item(91, 75).Value = 95
$ws1.Cells.Item(91, 76).Value = 1
$ws1.Cells.Item(91, 77).Value = 2
$ws1.Cells.Item(91, 78).Value = 4
$ws1.Cells.Item(91, 79).Value = 2
$ws1.Cells.Item(91, 80).Value = 3
$ws1.Cells.Item(91, 81).Value = 3
$ws1.Cells.Item(91, 82).Value = 0
$ws1.Cells.Item(91, 83).Value = 3
$ws1.Cells.Item(91, 84).Value = 0
$ws1.Cells.Item(91, 85).Value = 1
$ws1.Cells.Item(91, 86).Value = 15
$ws1.Cells.Item(91, 87).Value = 65
$ws1.Cells.Item(91, 88).Value = 73
$ws1.Cells.Item(91, 89).Value = 106
$ws1.Cells.Item(91, 90).Value = 109
$ws1.Cells.Item(91, 91).Value = 3076
$ws1.Cells.Item(91, 92).Value = 48
$ws1.Cells.Item(91, 93).Value = 294
$ws1.Cells.Item(91, 94).Value = 3
$ws1.Cells.Item(91, 95).Value = 584
$ws1.Cells.Item(91, 96).Value = 5974
$ws1.Cells.Item(91, 97).Value = 6
$ws1.Cells.Item(91, 98).Value = 194
$ws1.Cells.Item(91, 99).Value = 9
$ws1.Cells.Item(91, 100).Value = 23
$ws1.Cells.Item(91, 101).Value = 36
$ws1.Cells.Item(91, 102).Value = 7
$ws1.Cells.Item(91, 103).Value = 0
$ws1.Cells.Item(91, 104).Value = 10
$ws1.Cells.Item(91, 105).Value = 26
$ws1.Cells.Item(91, 106).Value = 0
$ws1.Cells.Item(91, 107).Value = 25
$ws1.Cells.Item(91, 108).Value = 21
$ws1.Cells.Item(91, 109).Value = 0
$ws1.Cells.Item(91, 110).Value = 0
$ws1.Cells.Item(91, 111).Value = 1
$ws1.Cells.Item(91, 112).Value = 36
$ws1.Cells.Item(91, 113).Value = 226
$ws1.Cells.Item(91, 114).Value = 4
$ws1.Cells.Item(91, 115).Value = 159
$ws1.Cells.Item(91, 116).Value = 2
$ws1.Cells.Item(91, 117).Value = 5
$ws1.Cells.Item(91, 118).Value = 0
$ws1.Cells.Item(91, 119).Value = 2
$ws1.Cells.Item(91, 120).Value = 8
$ws1.Cells.Item(91, 121).Value = 131
$ws1.Cells.Item(91, 122).Value = 32
$ws1.Cells.Item(91, 123).Value = 1
$ws1.Cells.Item(91, 124).Value = 1
$ws1.Cells.Item(91, 125).Value = 9
$ws1.Cells.Item(91, 126).Value = 0
$ws1.Cells.Item(91, 127).Value = 13
$ws1.Cells.Item(91, 128).Value = 1
$ws1.Cells.Item(91, 129).Value = 3
$ws1.Cells.Item(91, 130).Value = 0
$ws1.Cells.Item(91, 131).Value = 2
$ws1.Cells.Item(91, 132).Value = 1159
$ws1.Cells.Item(91, 133).Value = 0
$ws1.Cells.Item(91, 134).Value = 60
$ws1.Cells.Item(91, 135).Value = 1
$ws1.Cells.Item(91, 136).Value = 0
$ws1.Cells.Item(91, 137).Value = 16
$ws1.Cells.Item(91, 138).Value = 22
$ws1.Cells.Item(91, 139).Value = 264
$ws1.Cells.Item(91, 140).Value = 21
$ws1.Cells.Item(91, 141).Value = 211
$ws1.Cells.Item(91, 142).Value = 6
$ws1.Cells.Item(91, 143).Value = 115
$ws1.Cells.Item(91, 144).Value = 0
$ws1.Cells.Item(91, 145).Value = 3
$ws1.Cells.Item(91, 146).Value = 36
$ws1.Cells.Item(91, 147).Value = 272
$ws1.Cells.Item(91, 148).Value = 249
$ws1.Cells.Item(91, 149).Value = 902
$ws1.Cells.Item(91, 150).Value = 36
$ws1.Cells.Item(91, 151).Value = 28
$ws1.Cells.Item(91, 152).Value = 160
$ws1.Cells.Item(91, 153).Value = 228
$ws1.Cells.Item(91, 154).Value = 6
$ws1.Cells.Item(91, 155).Value = 0
$ws1.Cells.Item(91, 156).Value = 1
$ws1.Cells.Item(91, 157).Value = 0
$ws1.Cells.Item(91, 158).Value = 1
$ws1.Cells.Item(91, 159).Value = 99
$ws1.Cells.Item(91, 160).Value = 11
$ws1.Cells.Item(91, 161).Value = 202
$ws1.Cells.Item(91, 162).Value = 0
$ws1.Cells.Item(91, 163).Value = 71
$ws1.Cells.Item(91, 164).Value = 0
$ws1.Cells.Item(91, 165).Value = 0
$ws1.Cells.Item(91, 166).Value = 59
$ws1.Cells.Item(91, 167).Value = 0
$ws1.Cells.Item(91, 168).Value = 17
$ws1.Cells.Item(91, 169).Value = 105
$ws1.Cells.Item(91, 170).Value = 8189
$ws1.Cells.Item(91, 171).Value = 9
$ws1.Cells.Item(91, 172).Value = 2
$ws1.Cells.Item(91, 173).Value = 0
$ws1.Cells.Item(91, 174).Value = 401
$ws1.Cells.Item(91, 175).Value = 1048
$ws1.Cells.Item(91, 176).Value = 0
$ws1.Cells.Item(91, 177).Value = 16
$ws1.Cells.Item(91, 178).Value = 109
$ws1.Cells.Item(91, 179).Value = 0
$ws1.Cells.Item(91, 180).Value = 3
$ws1.Cells.Item(91, 181).Value = 8
$ws1.Cells.Item(91, 182).Value = 0
$ws1.Cells.Item(91, 183).Value = 1704
$ws1.Cells.Item(91, 184).Value = 3
$ws1.Cells.Item(91, 185).Value = 16
$ws1.Cells.Item(91, 186).Value = 93
$ws1.Cells.Item(91, 187).Value = 63
$ws1.Cells.Item(91, 188).Value = 2546
$ws1.Cells.Item(91, 189).Value = 0
$ws1.Cells.Item(91, 190).Value = 3
$ws1.Cells.Item(91, 191).Value = 19979
$ws1.Cells.Item(91, 192).Value = 66
$ws1.Cells.Item(91, 193).Value = 29
$ws1.Cells.Item(91, 194).Value = 0
$ws1.Cells.Item(91, 195).Value = 54
$ws1.Cells.Item(91, 196).Value = 12
$ws1.Cells.Item(91, 197).Value = 2

# --- Deaths sheet (sheet2), row 91 ---
$ws2.Cells.Item(91, 2).Value = 1
$ws2.Cells.Item(91, 3).Value = 1
$ws2.Cells.Item(91, 4).Value = 5
$ws2.Cells.Item(91, 5).Value = 1
$ws2.Cells.Item(91, 6).Value = 0
$ws2.Cells.Item(91, 7).Value = 0
$ws2.Cells.Item(91, 8).Value = 0
$ws2.Cells.Item(91, 9).Value = 2
$ws2.Cells.Item(91, 10).Value = 2
$ws2.Cells.Item(91, 11).Value = 0
$ws2.Cells.Item(91, 12).Value = 1
$ws2.Cells.Item(91, 13).Value = 0
$ws2.Cells.Item(91, 14).Value = 1
$ws2.Cells.Item(91, 15).Value = 0
$ws2.Cells.Item(91, 16).Value = 0
$ws2.Cells.Item(91, 17).Value = 0
$ws2.Cells.Item(91, 18).Value = 0
$ws2.Cells.Item(91, 19).Value = 0
$ws2.Cells.Item(91, 20).Value = 64
$ws2.Cells.Item(91, 21).Value = 0
$ws2.Cells.Item(91, 22).Value = 0
$ws2.Cells.Item(91, 23).Value = 0
$ws2.Cells.Item(91, 24).Value = 0
$ws2.Cells.Item(91, 25).Value = 0
$ws2.Cells.Item(91, 26).Value = 2
$ws2.Cells.Item(91, 27).Value = 22
$ws2.Cells.Item(91, 28).Value = 0
$ws2.Cells.Item(91, 29).Value = 1
$ws2.Cells.Item(91, 30).Value = 4
$ws2.Cells.Item(91, 31).Value = 6
$ws2.Cells.Item(91, 32).Value = 0
$ws2.Cells.Item(91, 33).Value = 1
$ws2.Cells.Item(91, 34).Value = 7
$ws2.Cells.Item(91, 35).Value = 0
$ws2.Cells.Item(91, 37).Value = 0
$ws2.Cells.Item(91, 38).Value = 0
$ws2.Cells.Item(91, 39).Value = 0
$ws2.Cells.Item(91, 40).Value = 1
$ws2.Cells.Item(91, 41).Value = 5
$ws2.Cells.Item(91, 42).Value = 0
$ws2.Cells.Item(91, 43).Value = 0
$ws2.Cells.Item(91, 44).Value = 0
$ws2.Cells.Item(91, 45).Value = 0
$ws2.Cells.Item(91, 46).Value = 2
$ws2.Cells.Item(91, 47).Value = 1
$ws2.Cells.Item(91, 48).Value = 0
$ws2.Cells.Item(91, 49).Value = 0
$ws2.Cells.Item(91, 50).Value = 2
$ws2.Cells.Item(91, 51).Value = 2
$ws2.Cells.Item(91, 52).Value = 13
$ws2.Cells.Item(91, 53).Value = 0
$ws2.Cells.Item(91, 54).Value = 0
$ws2.Cells.Item(91, 55).Value = 8
$ws2.Cells.Item(91, 56).Value = 7
$ws2.Cells.Item(91, 57).Value = 6
$ws2.Cells.Item(91, 58).Value = 0
$ws2.Cells.Item(91, 59).Value = 0
$ws2.Cells.Item(91, 60).Value = 0
$ws2.Cells.Item(91, 61).Value = 0
$ws2.Cells.Item(91, 62).Value = 0
$ws2.Cells.Item(91, 63).Value = 0
$ws2.Cells.Item(91, 64).Value = 0
$ws2.Cells.Item(91, 65).Value = 0
$ws2.Cells.Item(91, 66).Value = 2
$ws2.Cells.Item(91, 67).Value = 319
$ws2.Cells.Item(91, 68).Value = 0
$ws2.Cells.Item(91, 69).Value = 0
$ws2.Cells.Item(91, 70).Value = 0
$ws2.Cells.Item(91, 71).Value = 0
$ws2.Cells.Item(91, 72).Value = 64
$ws2.Cells.Item(91, 73).Value = 2
$ws2.Cells.Item(91, 74).Value = 0
$ws2.Cells.Item(91, 75).Value = 4
$ws2.Cells.Item(91, 76).Value = 0
$ws2.Cells.Item(91, 77).Value = 0
$ws2.Cells.Item(91, 78).Value = 0
$ws2.Cells.Item(91, 79).Value = 0
$ws2.Cells.Item(91, 80).Value = 0
$ws2.Cells.Item(91, 81).Value = 0
$ws2.Cells.Item(91, 82).Value = 0
$ws2.Cells.Item(91, 83).Value = 0
$ws2.Cells.Item(91, 84).Value = 0
$ws2.Cells.Item(91, 85).Value = 0
$ws2.Cells.Item(91, 86).Value = 1
$ws2.Cells.Item(91, 87).Value = 2
$ws2.Cells.Item(91, 88).Value = 0
$ws2.Cells.Item(91, 89).Value = 6
$ws2.Cells.Item(91, 90).Value = 15
$ws2.Cells.Item(91, 91).Value = 139
$ws2.Cells.Item(91, 92).Value = 2
$ws2.Cells.Item(91, 93).Value = 14
$ws2.Cells.Item(91, 94).Value = 0
$ws2.Cells.Item(91, 95).Value = 2
$ws2.Cells.Item(91, 96).Value = 887
$ws2.Cells.Item(91, 97).Value = 0
$ws2.Cells.Item(91, 98).Value = 3
$ws2.Cells.Item(91, 99).Value = 0
$ws2.Cells.Item(91, 100).Value = 1
$ws2.Cells.Item(91, 101).Value = 0
$ws2.Cells.Item(91, 102).Value = 0
$ws2.Cells.Item(91, 103).Value = 0
$ws2.Cells.Item(91, 104).Value = 0
$ws2.Cells.Item(91, 105).Value = 0
$ws2.Cells.Item(91, 106).Value = 0
$ws2.Cells.Item(91, 107).Value = 0
$ws2.Cells.Item(91, 108).Value = 1
$ws2.Cells.Item(91, 109).Value = 0
$ws2.Cells.Item(91, 110).Value = 0
$ws2.Cells.Item(91, 111).Value = 0
$ws2.Cells.Item(91, 112).Value = 2
$ws2.Cells.Item(91, 113).Value = 3
$ws2.Cells.Item(91, 114).Value = 0
$ws2.Cells.Item(91, 115).Value = 1
$ws2.Cells.Item(91, 116).Value = 0
$ws2.Cells.Item(91, 117).Value = 0
$ws2.Cells.Item(91, 118).Value = 0
$ws2.Cells.Item(91, 119).Value = 0
$ws2.Cells.Item(91, 120).Value = 0
$ws2.Cells.Item(91, 121).Value = 4
$ws2.Cells.Item(91, 122).Value = 0
$ws2.Cells.Item(91, 123).Value = 0
$ws2.Cells.Item(91, 124).Value = 0
$ws2.Cells.Item(91, 125).Value = 0
$ws2.Cells.Item(91, 126).Value = 0
$ws2.Cells.Item(91, 127).Value = 0
$ws2.Cells.Item(91, 128).Value = 0
$ws2.Cells.Item(91, 129).Value = 0
$ws2.Cells.Item(91, 130).Value = 0
$ws2.Cells.Item(91, 131).Value = 0
$ws2.Cells.Item(91, 132).Value = 93
$ws2.Cells.Item(91, 133).Value = 0
$ws2.Cells.Item(91, 134).Value = 1
$ws2.Cells.Item(91, 135).Value = 0
$ws2.Cells.Item(91, 136).Value = 0
$ws2.Cells.Item(91, 137).Value = 0
$ws2.Cells.Item(91, 138).Value = 1
$ws2.Cells.Item(91, 139).Value = 4
$ws2.Cells.Item(91, 140).Value = 0
$ws2.Cells.Item(91, 141).Value = 2
$ws2.Cells.Item(91, 142).Value = 0
$ws2.Cells.Item(91, 143).Value = 3
$ws2.Cells.Item(91, 144).Value = 0
$ws2.Cells.Item(91, 145).Value = 0
$ws2.Cells.Item(91, 146).Value = 5
$ws2.Cells.Item(91, 147).Value = 14
$ws2.Cells.Item(91, 148).Value = 2
$ws2.Cells.Item(91, 149).Value = 24
$ws2.Cells.Item(91, 150).Value = 1
$ws2.Cells.Item(91, 151).Value = 1
$ws2.Cells.Item(91, 152).Value = 5
$ws2.Cells.Item(91, 153).Value = 1
$ws2.Cells.Item(91, 154).Value = 0
$ws2.Cells.Item(91, 155).Value = 0
$ws2.Cells.Item(91, 156).Value = 0
$ws2.Cells.Item(91, 157).Value = 0
$ws2.Cells.Item(91, 158).Value = 1
$ws2.Cells.Item(91, 159).Value = 1
$ws2.Cells.Item(91, 160).Value = 0
$ws2.Cells.Item(91, 161).Value = 4
$ws2.Cells.Item(91, 162).Value = 0
$ws2.Cells.Item(91, 163).Value = 1
$ws2.Cells.Item(91, 164).Value = 0
$ws2.Cells.Item(91, 165).Value = 0
$ws2.Cells.Item(91, 166).Value = 0
$ws2.Cells.Item(91, 167).Value = 0
$ws2.Cells.Item(91, 168).Value = 0
$ws2.Cells.Item(91, 169).Value = 8
$ws2.Cells.Item(91, 170).Value = 832
$ws2.Cells.Item(91, 171).Value = 1
$ws2.Cells.Item(91, 172).Value = 0
$ws2.Cells.Item(91, 173).Value = 0
$ws2.Cells.Item(91, 174).Value = 10
$ws2.Cells.Item(91, 175).Value = 38
$ws2.Cells.Item(91, 176).Value = 0
$ws2.Cells.Item(91, 177).Value = 0
$ws2.Cells.Item(91, 178).Value = 1
$ws2.Cells.Item(91, 179).Value = 0
$ws2.Cells.Item(91, 180).Value = 0
$ws2.Cells.Item(91, 181).Value = 0
$ws2.Cells.Item(91, 182).Value = 0
$ws2.Cells.Item(91, 183).Value = 16
$ws2.Cells.Item(91, 184).Value = 0
$ws2.Cells.Item(91, 185).Value = 0
$ws2.Cells.Item(91, 186).Value = 3
$ws2.Cells.Item(91, 187).Value = 0
$ws2.Cells.Item(91, 188).Value = 260
$ws2.Cells.Item(91, 189).Value = 1
$ws2.Cells.Item(91, 190).Value = 0
$ws2.Cells.Item(91, 191).Value = 484
$ws2.Cells.Item(91, 192).Value = 0
$ws2.Cells.Item(91, 193).Value = 1
$ws2.Cells.Item(91, 194).Value = 0
$ws2.Cells.Item(91, 195).Value = 0
$ws2.Cells.Item(91, 196).Value = 0
$ws2.Cells.Item(91, 197).Value = 0
